# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
# Updates the "Croatia HNL" sheet with newly scraped match rows.
#
# Net effect (per the target diff):
#   - Existing rows 146-148 get overwritten with "next" match data
#     (each row's content effectively shifts down by one row vs. the
#     previous snapshot, and a brand-new match is inserted at the top
#     of this block).
#   - Two brand-new rows (149, 150) are appended at the bottom, reusing
#     the same look (bold/centered/bordered index column, date-formatted
#     Date column) as all the other data rows.
#   - The sheet's used range grows from A1:AC148 to A1:AC150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 149-150 are brand new - clone the formatting of an existing data
# row (148) onto them first, so the "id index" column (A) keeps its
# bold/centered/bordered look and the "Date" column (E) keeps its
# custom date-time number format, matching every other row on the sheet.
$ws.Range("A148:G148").Copy()
$ws.Range("A149:G150").PasteSpecial(-4122)
$ws.Range("K148:AA148").Copy()
$ws.Range("K149:AA150").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

function Set-DataRow {
    param(
        [int]$Row,
        [int]$Idx,
        [int]$Id,
        [double]$DateSerial,
        [string]$HomeTeam,
        [string]$AwayTeam,
        [double]$K, [double]$L, [double]$M,
        [double]$N, [double]$O, [double]$P, [double]$Q,
        [double]$R, [double]$S, [double]$T, [double]$U, [double]$V
    )

    $ws.Cells.Item($Row, 1).Value = $Idx                # A: running index
    $ws.Cells.Item($Row, 2).Value = $Id                 # B: id
    $ws.Cells.Item($Row, 3).Value = "Croatia HNL"        # C: Div
    $ws.Cells.Item($Row, 4).Value = "Croatia HNL"        # D: Div Original Name
    $ws.Cells.Item($Row, 5).Value = $DateSerial          # E: Date
    $ws.Cells.Item($Row, 6).Value = $HomeTeam            # F: HomeTeam
    $ws.Cells.Item($Row, 7).Value = $AwayTeam            # G: AwayTeam

    $ws.Cells.Item($Row, 11).Value = $K                  # K: oddH_op
    $ws.Cells.Item($Row, 12).Value = $L                  # L: oddD_op
    $ws.Cells.Item($Row, 13).Value = $M                  # M: oddA_op
    $ws.Cells.Item($Row, 14).Value = $N                  # N: oddH
    $ws.Cells.Item($Row, 15).Value = $O                  # O: oddD
    $ws.Cells.Item($Row, 16).Value = $P                  # P: oddA
    $ws.Cells.Item($Row, 17).Value = $Q                  # Q: Ah

    $ws.Cells.Item($Row, 18).Value = $R                  # R: oddAHH
    $ws.Cells.Item($Row, 19).Value = $S                  # S: oddAHA
    $ws.Cells.Item($Row, 20).Value = $T                  # T: AhOU
    $ws.Cells.Item($Row, 21).Value = $U                  # U: oddAHOver
    $ws.Cells.Item($Row, 22).Value = $V                  # V: oddAHUnder

    $ws.Cells.Item($Row, 23).Value = 0                   # W: PLH
    $ws.Cells.Item($Row, 24).Value = 0                   # X: PLD
    $ws.Cells.Item($Row, 25).Value = 0                   # Y: PLA
    $ws.Cells.Item($Row, 26).Value = 0                   # Z: PL_Ahh
    $ws.Cells.Item($Row, 27).Value = 0                   # AA: PL_Aha
}

# NOTE: this runtime's PowerShell subset only binds POSITIONAL arguments
# reliably, so Set-DataRow is always called without `-ParamName` syntax.

# Row 146 — NK Rudes vs Slaven Belupo (new top match of this block)
Set-DataRow 146 144 6769308 45394.54166666666 `
    "NK Rudes" "Slaven Belupo" `
    4 3.4 1.95 4.5 3.5 1.833 0.5 `
    2.025 1.825 2.5 1.975 1.875

# Row 147 — NK Varazdin vs NK Lokomotiva Zagreb
Set-DataRow 147 145 6788942 45395.49305555555 `
    "NK Varazdin" "NK Lokomotiva Zagreb" `
    2.8 3.25 2.5 2.8 3.25 2.5 0 `
    2.025 1.825 2.5 1.975 1.875

# Row 148 — HNK Gorica vs Dinamo Zagreb
Set-DataRow 148 146 6788943 45395.58333333334 `
    "HNK Gorica" "Dinamo Zagreb" `
    8 4.5 1.363 8.5 4.5 1.333 1.5 `
    1.825 2.025 2.75 1.925 1.925

# Row 149 — Hajduk Split vs NK Osijek (brand-new row)
Set-DataRow 149 147 6923266 45396.47916666666 `
    "Hajduk Split" "NK Osijek" `
    1.615 3.5 6 1.6 3.5 6 -1 `
    2.1 1.775 2.25 1.825 2.025

# Row 150 — Istra 1961 vs HNK Rijeka (brand-new row)
Set-DataRow 150 148 6788944 45396.58333333334 `
    "Istra 1961" "HNK Rijeka" `
    5.5 3.6 1.615 5.5 3.6 1.615 0.75 `
    2.025 1.825 2.25 1.875 1.975
